$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.305.81"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "1.857.40"
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("E4").Value = "  -0.73%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.00"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4621"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3702"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07317"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8835"
$ws.Range("D10").ClearFormats()
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.82"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07807"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").Value = "1.837.33"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.544"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.85"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008868"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.80"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("D21").Value = "27.328.13"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.123"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("E23").Value = "  -1.01%  "
$ws.Range("D24").Value = "2.057.42"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.909"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.18"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.075"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.117"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.94"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08860"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7628"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.52%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.178"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.87%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.994"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.497"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.608"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.16%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01963"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.076"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.989"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05210"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.028"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5152"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1640"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.353"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4838"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.32"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.13"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.653"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.56"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.76%  "
